$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new daily-log row at the bottom of the sheet:
#   A46 = 2025/10/01 (date, stored as text like every other row above)
#   B46 = 水         (weekday, text)
#   C46 = 20         (hour, number)
#   D46 = 19         (ranking, number)
#
# The date/weekday columns in this sheet are plain text, not real Excel
# dates (see rows 2-45). Assigning "2025/10/01" directly would be
# auto-recognized as a date literal and stored as a serial number, so a
# leading apostrophe forces text entry instead. ClearFormats() afterwards
# drops the resulting "quote prefix" cell style so the new row carries no
# formatting beyond what the existing rows already have.
$ws.Range("A46").Value = "'2025/10/01"
$ws.Range("B46").Value = "水"
$ws.Range("C46").Value = 20
$ws.Range("D46").Value = 19
$ws.Range("A46").ClearFormats()
